$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5396133333333334
$ws.Range("H2").Value = 1.61884
$ws.Range("I2").Value = 0.02587018426425635
$ws.Range("J2").Value = 0.02587018426425635
$ws.Range("M2").Value = 0.5396133333333334
$ws.Range("N2").Value = 1.61884
$ws.Range("O2").Value = 0.02587018426425635
$ws.Range("P2").Value = 0.02587018426425635
$ws.Range("Q2").Value = 0.2911825495111112
$ws.Range("R2").Value = 2.6206429456
$ws.Range("S2").Value = 0.0006692664338665771
$ws.Range("T2").Value = 0.0006692664338665771

# Row 3
$ws.Range("G3").Value = 0.5396133333333334
$ws.Range("H3").Value = 1.61884
$ws.Range("I3").Value = 0.02587018426425635
$ws.Range("J3").Value = 0.02587018426425635
$ws.Range("M3").Value = 6.260434
$ws.Range("N3").Value = 18.781302
$ws.Range("O3").Value = 0.3001382122153186
$ws.Range("P3").Value = 0.3001382122153186
$ws.Range("Q3").Value = 3.378213658853334
$ws.Range("R3").Value = 30.40392292968
$ws.Range("S3").Value = 0.007764630854754769
$ws.Range("T3").Value = 0.007764630854754769

# Row 4
$ws.Range("G4").Value = 0.5396133333333334
$ws.Range("H4").Value = 1.61884
$ws.Range("I4").Value = 0.02587018426425635
$ws.Range("J4").Value = 0.02587018426425635
$ws.Range("M4").Value = 1.182122
$ws.Range("N4").Value = 3.546366
$ws.Range("O4").Value = 0.0566733845769154
$ws.Range("P4").Value = 0.0566733845769154
$ws.Range("Q4").Value = 0.6378887928266667
$ws.Range("R4").Value = 5.74099913544
$ws.Range("S4").Value = 0.001466150901883866
$ws.Range("T4").Value = 0.001466150901883866

# Row 5
$ws.Range("G5").Value = 0.5396133333333334
$ws.Range("H5").Value = 1.61884
$ws.Range("I5").Value = 0.02587018426425635
$ws.Range("J5").Value = 0.02587018426425635
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.87633433333333
$ws.Range("N5").Value = 38.629003
$ws.Range("O5").Value = 0.6173182189435098
$ws.Range("P5").Value = 0.6173182189435098
$ws.Range("Q5").Value = 6.948241690724445
$ws.Range("R5").Value = 62.53417521652
$ws.Range("S5").Value = 0.01597013607375115
$ws.Range("T5").Value = 0.01597013607375115

# Row 6
$ws.Range("G6").Value = 6.260434
$ws.Range("H6").Value = 18.781302
$ws.Range("I6").Value = 0.3001382122153186
$ws.Range("J6").Value = 0.3001382122153186
$ws.Range("M6").Value = 0.5396133333333334
$ws.Range("N6").Value = 1.61884
$ws.Range("O6").Value = 0.02587018426425635
$ws.Range("P6").Value = 0.02587018426425635
$ws.Range("Q6").Value = 3.378213658853334
$ws.Range("R6").Value = 30.40392292968
$ws.Range("S6").Value = 0.007764630854754769
$ws.Range("T6").Value = 0.007764630854754769

# Row 7
$ws.Range("G7").Value = 6.260434
$ws.Range("H7").Value = 18.781302
$ws.Range("I7").Value = 0.3001382122153186
$ws.Range("J7").Value = 0.3001382122153186
$ws.Range("M7").Value = 6.260434
$ws.Range("N7").Value = 18.781302
$ws.Range("O7").Value = 0.3001382122153186
$ws.Range("P7").Value = 0.3001382122153186
$ws.Range("Q7").Value = 39.193033868356
$ws.Range("R7").Value = 352.737304815204
$ws.Range("S7").Value = 0.09008294643180763
$ws.Range("T7").Value = 0.09008294643180763

# Row 8
$ws.Range("G8").Value = 6.260434
$ws.Range("H8").Value = 18.781302
$ws.Range("I8").Value = 0.3001382122153186
$ws.Range("J8").Value = 0.3001382122153186
$ws.Range("M8").Value = 1.182122
$ws.Range("N8").Value = 3.546366
$ws.Range("O8").Value = 0.0566733845769154
$ws.Range("P8").Value = 0.0566733845769154
$ws.Range("Q8").Value = 7.400596760947999
$ws.Range("R8").Value = 66.605370848532
$ws.Range("S8").Value = 0.0170098483271066
$ws.Range("T8").Value = 0.0170098483271066

# Row 9
$ws.Range("G9").Value = 6.260434
$ws.Range("H9").Value = 18.781302
$ws.Range("I9").Value = 0.3001382122153186
$ws.Range("J9").Value = 0.3001382122153186
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 12.87633433333333
$ws.Range("N9").Value = 38.629003
$ws.Range("O9").Value = 0.6173182189435098
$ws.Range("P9").Value = 0.6173182189435098
$ws.Range("Q9").Value = 80.61144125576733
$ws.Range("R9").Value = 725.502971301906
$ws.Range("S9").Value = 0.1852807866016496
$ws.Range("T9").Value = 0.1852807866016496

# Row 10
$ws.Range("G10").Value = 1.182122
$ws.Range("H10").Value = 3.546366
$ws.Range("I10").Value = 0.0566733845769154
$ws.Range("J10").Value = 0.0566733845769154
$ws.Range("M10").Value = 0.5396133333333334
$ws.Range("N10").Value = 1.61884
$ws.Range("O10").Value = 0.02587018426425635
$ws.Range("P10").Value = 0.02587018426425635
$ws.Range("Q10").Value = 0.6378887928266667
$ws.Range("R10").Value = 5.74099913544
$ws.Range("S10").Value = 0.001466150901883866
$ws.Range("T10").Value = 0.001466150901883866

# Row 11
$ws.Range("G11").Value = 1.182122
$ws.Range("H11").Value = 3.546366
$ws.Range("I11").Value = 0.0566733845769154
$ws.Range("J11").Value = 0.0566733845769154
$ws.Range("M11").Value = 6.260434
$ws.Range("N11").Value = 18.781302
$ws.Range("O11").Value = 0.3001382122153186
$ws.Range("P11").Value = 0.3001382122153186
$ws.Range("Q11").Value = 7.400596760947999
$ws.Range("R11").Value = 66.605370848532
$ws.Range("S11").Value = 0.0170098483271066
$ws.Range("T11").Value = 0.0170098483271066

# Row 12
$ws.Range("G12").Value = 1.182122
$ws.Range("H12").Value = 3.546366
$ws.Range("I12").Value = 0.0566733845769154
$ws.Range("J12").Value = 0.0566733845769154
$ws.Range("M12").Value = 1.182122
$ws.Range("N12").Value = 3.546366
$ws.Range("O12").Value = 0.0566733845769154
$ws.Range("P12").Value = 0.0566733845769154
$ws.Range("Q12").Value = 1.397412422884
$ws.Range("R12").Value = 12.576711805956
$ws.Range("S12").Value = 0.003211872519402953
$ws.Range("T12").Value = 0.003211872519402953

# Row 13
$ws.Range("G13").Value = 1.182122
$ws.Range("H13").Value = 3.546366
$ws.Range("I13").Value = 0.0566733845769154
$ws.Range("J13").Value = 0.0566733845769154
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 12.87633433333333
$ws.Range("N13").Value = 38.629003
$ws.Range("O13").Value = 0.6173182189435098
$ws.Range("P13").Value = 0.6173182189435098
$ws.Range("Q13").Value = 15.22139809478866
$ws.Range("R13").Value = 136.992582853098
$ws.Range("S13").Value = 0.034985512828522
$ws.Range("T13").Value = 0.034985512828522

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 12.87633433333333
$ws.Range("H14").Value = 38.629003
$ws.Range("I14").Value = 0.6173182189435098
$ws.Range("J14").Value = 0.6173182189435098
$ws.Range("M14").Value = 0.5396133333333334
$ws.Range("N14").Value = 1.61884
$ws.Range("O14").Value = 0.02587018426425635
$ws.Range("P14").Value = 0.02587018426425635
$ws.Range("Q14").Value = 6.948241690724445
$ws.Range("R14").Value = 62.53417521652
$ws.Range("S14").Value = 0.01597013607375115
$ws.Range("T14").Value = 0.01597013607375115

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 12.87633433333333
$ws.Range("H15").Value = 38.629003
$ws.Range("I15").Value = 0.6173182189435098
$ws.Range("J15").Value = 0.6173182189435098
$ws.Range("M15").Value = 6.260434
$ws.Range("N15").Value = 18.781302
$ws.Range("O15").Value = 0.3001382122153186
$ws.Range("P15").Value = 0.3001382122153186
$ws.Range("Q15").Value = 80.61144125576733
$ws.Range("R15").Value = 725.502971301906
$ws.Range("S15").Value = 0.1852807866016496
$ws.Range("T15").Value = 0.1852807866016496

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 12.87633433333333
$ws.Range("H16").Value = 38.629003
$ws.Range("I16").Value = 0.6173182189435098
$ws.Range("J16").Value = 0.6173182189435098
$ws.Range("M16").Value = 1.182122
$ws.Range("N16").Value = 3.546366
$ws.Range("O16").Value = 0.0566733845769154
$ws.Range("P16").Value = 0.0566733845769154
$ws.Range("Q16").Value = 15.22139809478866
$ws.Range("R16").Value = 136.992582853098
$ws.Range("S16").Value = 0.034985512828522
$ws.Range("T16").Value = 0.034985512828522

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 12.87633433333333
$ws.Range("H17").Value = 38.629003
$ws.Range("I17").Value = 0.6173182189435098
$ws.Range("J17").Value = 0.6173182189435098
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 12.87633433333333
$ws.Range("N17").Value = 38.629003
$ws.Range("O17").Value = 0.6173182189435098
$ws.Range("P17").Value = 0.6173182189435098
$ws.Range("Q17").Value = 165.7999858637787
$ws.Range("R17").Value = 1492.199872774009
$ws.Range("S17").Value = 0.3810817834395871
$ws.Range("T17").Value = 0.3810817834395871

